$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row "Differential Equations and Boundary Value Problems - Computing and
# Modeling - C. Henry Edwards, David E. Penney, David T. Calvis (2015,
# Pearson) 5th Edition.pdf" (old row 17) was removed from the book list.
# The S.No. column (A) keeps its original sequential numbering per row, but
# the FileName/Book/Author(s)/Link/Edn/Year/Publisher columns (B:H) for every
# subsequent row shift up by one row, and the final (now empty) row 45 is
# deleted.
$ws.Range("B17:H44").Value2 = $ws.Range("B18:H45").Value2
$ws.Rows(45).Delete() | Out-Null
